# The ExcelDataSourceFile sample data (previously living in A1:C3) is moved
# to F1:H3 on the same worksheet, preserving cell contents/order.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$source = $ws.Range("A1:C3")
$destination = $ws.Range("F1:H3")

# Copy the values over to their new home, then clear out the old range.
$destination.Value = $source.Value()
$source.Clear()

# Mirror the author's selection ending up on the relocated table.
[void]$destination.Select()
